# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> "Office Theme" (used by the notes master)
#   ppt/theme/theme2.xml  -> "Integral"     (used by the slide master / all slides)
# The target edit swaps the two themes' contents, so the slide master
# (and therefore every slide) ends up rendered with the "Office Theme"
# palette instead of "Integral".
#
# The font scheme and format scheme are identical between the two themes;
# only the 12 theme colors (and the cosmetic <a:theme>/<a:clrScheme> name
# attributes, which PowerPoint's automation surface does not expose for
# writing) differ. So we reproduce the visible effect of the swap by
# pushing the "Office Theme" color values onto the presentation's live
# color scheme, which is backed by the theme part driving the slide
# master (theme2.xml).

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.ColorScheme

# Index -> (scheme slot, target RGB as 0xBBGGRR for the COM `RGB` setter)
# Values below are the "Office Theme" colors (currently in theme1.xml):
#   dk1=000000 lt1=FFFFFF dk2=44546A lt2=E7E6E6
#   accent1=5B9BD5 accent2=ED7D31 accent3=A5A5A5 accent4=FFC000
#   accent5=4472C4 accent6=70AD47 hlink=0563C1 folHlink=954F72
$officeColors = @(
    0,          # 1  dk1      000000
    16777215,   # 2  lt1      FFFFFF
    6968388,    # 3  dk2      44546A
    15132391,   # 4  lt2      E7E6E6
    13998939,   # 5  accent1  5B9BD5
    3243501,    # 6  accent2  ED7D31
    10855845,   # 7  accent3  A5A5A5
    49407,      # 8  accent4  FFC000
    12874308,   # 9  accent5  4472C4
    4697456,    # 10 accent6  70AD47
    12673797,   # 11 hlink    0563C1
    7491477     # 12 folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $cs.Colors($i).RGB = $officeColors[$i - 1]
}
